$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 617.4
$ws.Range("J2").Value = 559.2
$ws.Range("L2").Value = 559.2
$ws.Range("N2").Value = -785.2

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 837.4545000000001
$ws.Range("I15").Value = 837.4545000000001
$ws.Range("K15").Value = 2512.3635
$ws.Range("M15").Value = -2343.3635

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 6820.8335
$ws.Range("I38").Value = 6185
$ws.Range("J38").Value = 10000
$ws.Range("K38").Value = 18555
$ws.Range("L38").Value = 30000
$ws.Range("M38").Value = -18183
$ws.Range("N38").Value = -30744

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 191.66667
$ws.Range("I39").Value = 191.66667
$ws.Range("K39").Value = 575.00001
$ws.Range("M39").Value = -279.00001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H44").Value = 10000
$ws.Range("J44").Value = 10000
$ws.Range("L44").Value = 10000
$ws.Range("N44").Value = -10924

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H52").Value = 100
$ws.Range("J52").Value = 100
$ws.Range("L52").Value = 300
$ws.Range("N52").Value = -620

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H115").Value = 895.44446
$ws.Range("I115").Value = 893.4
$ws.Range("J115").Value = 898
$ws.Range("K115").Value = 2680.2
$ws.Range("L115").Value = 2694
$ws.Range("M115").Value = -1113.2
$ws.Range("N115").Value = -5828

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 3256.5
$ws.Range("I116").Value = 3594.4443
$ws.Range("J116").Value = 2918.5557
$ws.Range("K116").Value = 3594.4443
$ws.Range("L116").Value = 2918.5557
$ws.Range("M116").Value = -152.4443000000001
$ws.Range("N116").Value = -9802.555700000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H118").Value = 1334.5714
$ws.Range("I118").Value = 1329.4
$ws.Range("K118").Value = 3988.2
$ws.Range("M118").Value = -2331.2

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 5749.375
$ws.Range("I125").Value = 4000
$ws.Range("K125").Value = 36000
$ws.Range("M125").Value = -33540

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 16029.134
$ws.Range("J129").Value = 33030.145
$ws.Range("L129").Value = 99090.435
$ws.Range("N129").Value = -109090.435

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 3716.5676
$ws.Range("I132").Value = 3891.3547
$ws.Range("J132").Value = 2813.5
$ws.Range("K132").Value = 11674.0641
$ws.Range("L132").Value = 8440.5
$ws.Range("M132").Value = -9144.0641
$ws.Range("N132").Value = -13500.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 3860.1667
$ws.Range("I137").Value = 4041.3333
$ws.Range("J137").Value = 3316.6667
$ws.Range("K137").Value = 12123.9999
$ws.Range("L137").Value = 9950.000100000001
$ws.Range("M137").Value = -9573.999899999999
$ws.Range("N137").Value = -15050.0001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 6177515.5
$ws.Range("J138").Value = 8552519
$ws.Range("L138").Value = 25657557
$ws.Range("N138").Value = -25667837

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1124.4783
$ws.Range("I2").Value = 835.5
$ws.Range("J2").Value = 1785
$ws.Range("K2").Value = 835.5
$ws.Range("L2").Value = 1785
$ws.Range("M2").Value = -722.5
$ws.Range("N2").Value = -2011

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11502253
$ws.Range("I32").Value = 15387738
$ws.Range("J32").Value = 22412.773
$ws.Range("K32").Value = 15387738
$ws.Range("L32").Value = 22412.773
$ws.Range("M32").Value = -15387451
$ws.Range("N32").Value = -22986.773

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 1124.4783
$ws.Range("I116").Value = 835.5
$ws.Range("J116").Value = 1785
$ws.Range("K116").Value = 835.5
$ws.Range("L116").Value = 1785
$ws.Range("M116").Value = 1458.5
$ws.Range("N116").Value = -6373

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H124").Value = 46058.285
$ws.Range("J124").Value = 46058.285
$ws.Range("L124").Value = 46058.285
$ws.Range("N124").Value = -55878.285

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H125").Value = 71500
$ws.Range("J125").Value = 71500
$ws.Range("L125").Value = 71500
$ws.Range("N125").Value = -81340

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H134").Value = 444998
$ws.Range("J134").Value = 444998
$ws.Range("L134").Value = 444998
$ws.Range("N134").Value = -455138

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1124.4783
$ws.Range("I3").Value = 835.5
$ws.Range("J3").Value = 1785
$ws.Range("K3").Value = 835.5
$ws.Range("L3").Value = 1785
$ws.Range("M3").Value = -721.5
$ws.Range("N3").Value = -2013

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 5378.9
$ws.Range("I134").Value = 5172.1333
$ws.Range("K134").Value = 15516.3999
$ws.Range("M134").Value = -12981.3999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 20839404
$ws.Range("I31").Value = 5749.75
$ws.Range("K31").Value = 5749.75
$ws.Range("M31").Value = -5454.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 20839404
$ws.Range("I34").Value = 5749.75
$ws.Range("K34").Value = 5749.75
$ws.Range("M34").Value = -5547.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 5287.3
$ws.Range("J58").Value = 5710.4287
$ws.Range("L58").Value = 5710.4287
$ws.Range("N58").Value = -6116.4287

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 4741.5884
$ws.Range("I99").Value = 4917.5835
$ws.Range("K99").Value = 4917.5835
$ws.Range("M99").Value = -3419.5835

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 4741.5884
$ws.Range("I126").Value = 4917.5835
$ws.Range("K126").Value = 14752.7505
$ws.Range("M126").Value = -12282.7505

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H131").Value = 76239
$ws.Range("J131").Value = 76239
$ws.Range("L131").Value = 76239
$ws.Range("N131").Value = -86319

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 76296.63
$ws.Range("I132").Value = 84999.21000000001
$ws.Range("K132").Value = 254997.63
$ws.Range("M132").Value = -252467.63

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 5287.3
$ws.Range("J136").Value = 5710.4287
$ws.Range("L136").Value = 17131.2861
$ws.Range("N136").Value = -22231.2861

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H141").Value = 291109.12
$ws.Range("J141").Value = 312497.75
$ws.Range("L141").Value = 312497.75
$ws.Range("N141").Value = -322857.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 10171
$ws.Range("I14").Value = 10171
$ws.Range("K14").Value = 30513
$ws.Range("M14").Value = -30340

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 842.875
$ws.Range("I121").Value = 688.8889
$ws.Range("J121").Value = 1040.8572
$ws.Range("K121").Value = 2066.6667
$ws.Range("L121").Value = 3122.5716
$ws.Range("M121").Value = -756.6667000000002
$ws.Range("N121").Value = -5742.571599999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 2384815
$ws.Range("J132").Value = 4449082.5
$ws.Range("L132").Value = 40041742.5
$ws.Range("N132").Value = -40046802.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H121").Value = 61999.332
$ws.Range("J121").Value = 61999.332
$ws.Range("L121").Value = 61999.332
$ws.Range("N121").Value = -65493.332

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 15796173
$ws.Range("I126").Value = 20012948
$ws.Range("K126").Value = 60038844
$ws.Range("M126").Value = -60036374

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2117.125
$ws.Range("I132").Value = 1702.4166
$ws.Range("K132").Value = 5107.2498
$ws.Range("M132").Value = -2577.2498

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2869.7273
$ws.Range("I100").Value = 2092.5557
$ws.Range("K100").Value = 2092.5557
$ws.Range("M100").Value = -1551.5557

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 3051.2766
$ws.Range("I126").Value = 3744.6765
$ws.Range("J126").Value = 1237.7693
$ws.Range("K126").Value = 11234.0295
$ws.Range("L126").Value = 3713.3079
$ws.Range("M126").Value = -8764.029500000001
$ws.Range("N126").Value = -8653.3079

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4921.827
$ws.Range("J132").Value = 3244.1428
$ws.Range("L132").Value = 9732.428400000001
$ws.Range("N132").Value = -14792.4284
